$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 72: 2025-09-05, 四方坪站 station
$ws.Range("A72").Value = 45905
$ws.Range("A72").NumberFormat = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'
$ws.Range("B72").Value = "四方坪站"
$ws.Range("C72").Value = 13055.83
$ws.Range("D72").Value = 10859.28
$ws.Range("E72").Value = 4573.97
$ws.Range("F72").Value = 518
$ws.Range("F72").NumberFormat = '0_);[Red]\(0\)'

# Row 73: 2025-09-05, 高岭站 station
$ws.Range("A73").Value = 45905
$ws.Range("A73").NumberFormat = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'
$ws.Range("B73").Value = "高岭站"
$ws.Range("C73").Value = 5577.43
$ws.Range("D73").Value = 4404.98
$ws.Range("E73").Value = 1487.14
$ws.Range("F73").Value = 203
$ws.Range("F73").NumberFormat = '0_);[Red]\(0\)'

$ws.Range("F78").Select()
